# PI-2023-24 Self-Assessment.xlsx - sprint 2 submission update
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: "Group and Self Assessment"
# ---------------------------------------------------------------------------
$wsGroup = $wb.Worksheets.Item("Group and Self Assessment")
$wsGroup.Range("D10").Value = 3

# ---------------------------------------------------------------------------
# Sheet: "User Stories"
# ---------------------------------------------------------------------------
$wsUS = $wb.Worksheets.Item("User Stories")

$wsUS.Range("B7").Value = 1230744
$wsUS.Range("C7").Value = 4

$wsUS.Range("B8").Value = 1230744
$wsUS.Range("C8").Value = 4

$wsUS.Range("B9").Value = 1230399
$wsUS.Range("D9").Value = 3

# newly filled-in user story rows 14-19
$wsUS.Range("A14").Value = 9
$wsUS.Range("B14").Value = 1230399
$wsUS.Range("C14").Value = 3
$wsUS.Range("D14").Value = 3

$wsUS.Range("A15").Value = 10
$wsUS.Range("B15").Value = 1191647
$wsUS.Range("C15").Value = 3
$wsUS.Range("D15").Value = 3

$wsUS.Range("A16").Value = 11
$wsUS.Range("B16").Value = 1230744
$wsUS.Range("C16").Value = 3
$wsUS.Range("D16").Value = 3

$wsUS.Range("A17").Value = 12
$wsUS.Range("B17").Value = 1230741
$wsUS.Range("C17").Value = 3
$wsUS.Range("D17").Value = 3

$wsUS.Range("A18").Value = 13
$wsUS.Range("B18").Value = 1191647
$wsUS.Range("C18").Value = 3
$wsUS.Range("D18").Value = 3

$wsUS.Range("A19").Value = 14
$wsUS.Range("B19").Value = 1231235
$wsUS.Range("C19").Value = 3
$wsUS.Range("D19").Value = 4

# update data validation ranges for the now-filled rows
$wsUS.Range("C18:C25").Validation.Delete()
$wsUS.Range("C6:C17").Validation.Delete()
$wsUS.Range("C20:C25").Validation.Add(3, 1, 3, "=`$E`$40:`$J`$40")
$wsUS.Range("C6:C19").Validation.Add(3, 1, 3, "=`$E`$3:`$J`$3")

# ---------------------------------------------------------------------------
# Sheet: "Project Development"
# ---------------------------------------------------------------------------
$wsPD = $wb.Worksheets.Item("Project Development")

$wsPD.Range("G5").Value = 3

$wsPD.Range("D6").Value = 4
$wsPD.Range("G6").Value = 3

$wsPD.Range("C7").Value = 3
$wsPD.Range("D7").Value = 3
$wsPD.Range("E7").Value = 3
$wsPD.Range("F7").Value = 4
$wsPD.Range("G7").Value = 3

# ---------------------------------------------------------------------------
# Sheet: "Project Management"
# ---------------------------------------------------------------------------
$wsPM = $wb.Worksheets.Item("Project Management")

$wsPM.Range("C12").Value = 3
$wsPM.Range("D12").Value = 4
$wsPM.Range("E12").Value = 4
$wsPM.Range("F12").Value = 4
$wsPM.Range("G12").Value = 3

$wsPM.Range("C13").Value = 4
$wsPM.Range("D13").Value = 4
$wsPM.Range("E13").Value = 4
$wsPM.Range("F13").Value = 4
$wsPM.Range("G13").Value = 3

$wsPM.Range("C14").Value = 3
$wsPM.Range("D14").Value = 4
$wsPM.Range("E14").Value = 4
$wsPM.Range("F14").Value = 4
$wsPM.Range("G14").Value = 3
